$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1964723900037332
$ws.Range("C2").Value = 4.95718404746367
$ws.Range("B3").Value = 0.2061449038404378
$ws.Range("C3").Value = 5.272444997864113
